$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# A leading single-quote on numeric-looking values forces Excel to
# store them as text, matching the source data (which is not truly numeric).
$updates = @(
    @{Cell="D2"; Value='67.012.95'},
    @{Cell="E2"; Value='  +0.24%  '},
    @{Cell="D3"; Value='3.123.05'},
    @{Cell="E3"; Value='  +1.10%  '},
    @{Cell="E4"; Value='  -0.03%  '},
    @{Cell="D5"; Value='''577.44'},
    @{Cell="E5"; Value='  -0.43%  '},
    @{Cell="D6"; Value='''173.16'},
    @{Cell="E6"; Value='  +2.61%  '},
    @{Cell="E7"; Value='  -0.04%  '},
    @{Cell="E8"; Value='  -0.63%  '},
    @{Cell="D9"; Value='''6.45'},
    @{Cell="E9"; Value='  -2.41%  '},
    @{Cell="E10"; Value='  -0.90%  '},
    @{Cell="E11"; Value='  +0.23%  '},
    @{Cell="E12"; Value='  -1.15%  '},
    @{Cell="D13"; Value='''37.14'},
    @{Cell="E13"; Value='  +1.85%  '},
    @{Cell="E14"; Value='  -1.24%  '},
    @{Cell="D15"; Value='3.639.52'},
    @{Cell="E15"; Value='  +1.11%  '},
    @{Cell="D16"; Value='66.958.00'},
    @{Cell="E16"; Value='  +0.12%  '},
    @{Cell="E17"; Value='  -0.33%  '},
    @{Cell="D18"; Value='3.119.17'},
    @{Cell="D19"; Value='''16.29'},
    @{Cell="E19"; Value='  +0.78%  '},
    @{Cell="D20"; Value='''476.32'},
    @{Cell="E20"; Value='  +2.22%  '},
    @{Cell="E21"; Value='  -0.48%  '},
    @{Cell="D22"; Value='''7.95'},
    @{Cell="E22"; Value='  +5.64%  '},
    @{Cell="D23"; Value='''84.02'},
    @{Cell="E23"; Value='  +0.72%  '},
    @{Cell="E24"; Value='  +2.61%  '},
    @{Cell="E25"; Value='  -2.89%  '},
    @{Cell="D26"; Value='''10.09'},
    @{Cell="E26"; Value='  -0.61%  '},
    @{Cell="E27"; Value='  +0.03%  '},
    @{Cell="E28"; Value='  -1.11%  '},
    @{Cell="E29"; Value='  -1.17%  '},
    @{Cell="E30"; Value='  +0.57%  '},
    @{Cell="D31"; Value='''28.58'},
    @{Cell="E31"; Value='  +1.13%  '},
    @{Cell="E32"; Value='  +0.61%  '},
    @{Cell="D33"; Value='0.0₃0950'},
    @{Cell="E33"; Value='  -7.39%  '},
    @{Cell="D34"; Value='''0.999'},
    @{Cell="E34"; Value='  -0.11%  '},
    @{Cell="D35"; Value='''5.86'},
    @{Cell="E35"; Value='  -0.51%  '},
    @{Cell="D36"; Value='''0.976'},
    @{Cell="E36"; Value='  -2.97%  '},
    @{Cell="D37"; Value='''47.07'},
    @{Cell="E37"; Value='  -0.08%  '},
    @{Cell="D38"; Value='''50.18'},
    @{Cell="E38"; Value='  -0.12%  '},
    @{Cell="E39"; Value='  -2.28%  '},
    @{Cell="E40"; Value='  -1.80%  '},
    @{Cell="E41"; Value='  +1.10%  '},
    @{Cell="E42"; Value='  -0.11%  '},
    @{Cell="D43"; Value='2.816.26'},
    @{Cell="E43"; Value='  +1.56%  '},
    @{Cell="D44"; Value='''383.12'},
    @{Cell="E44"; Value='  -0.37%  '},
    @{Cell="E45"; Value='  -1.69%  '},
    @{Cell="E46"; Value='  -9.81%  '},
    @{Cell="D47"; Value='''135.49'},
    @{Cell="E47"; Value='  +0.29%  '},
    @{Cell="E48"; Value='  -0.02%  '},
    @{Cell="E49"; Value='  +0.14%  '},
    @{Cell="E50"; Value='  -1.48%  '},
    @{Cell="E51"; Value='  -0.64%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
